$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.330.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.178.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.61"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.606"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.73"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.11%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.56"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.89%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.76"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.497.68"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.15"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.172.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.766"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.248.94"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.86"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.70"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.14"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.63%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.43"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.32%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.38"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.98"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0809"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.11"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.30%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.107"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.23"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0335"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.17%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.77"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.196"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "59.34"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.13"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.77"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.472"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +14.08%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.17%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0968"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.22"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.09"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.36%  "
